# Fix part_list to load all columns dynamically
#
# The "sources" column (BC) for the two newly-discovered parts (rows 4 and 5,
# part numbers NP004 "Steel Screw" and NP007 "Aluminum Washer") was missing
# some of the source records that the loader should have picked up once it
# started reading every column dynamically. Re-run of the loader appended the
# missing "user" diagram references into the JSON array, and the refreshed
# rows (plus everything re-saved alongside them) picked up new `updated_at`
# timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Backfill the missing source entries in the "sources" JSON column ---

$ws.Range("BC4").Value2 = '[{"source_system": "pos", "source_file": "po_mock_newparts.pdf"}, {"source_system": "user", "source_file": "diagram2.pdf"}, {"source": "user", "file": "diagram 7.pdf", "description": "Steel Screw 39x90 mm"}]'

$ws.Range("BC5").Value2 = '[{"source_system": "pos", "source_file": "po_mock_newparts.pdf"}, {"source_system": "user", "source_file": "diagram 3.pdf"}, {"source_system": "user", "source_file": "diagram 3.pdf"}, {"source_system": "user", "source_file": "diagram 3.pdf"}, {"source": "user", "file": "diagram 7.pdf", "description": "Aluminum Washer 19x142 mm"}]'

# --- 2. Refresh "updated_at" (column C) timestamps that moved as a result ---

# Rows 3-5 (NP001, NP004, NP007) each got their own refresh timestamp,
# captured moments apart while the loader processed the new-parts rows.
$ws.Range("C3").Value2 = 46001.56353754582
$ws.Range("C4").Value2 = 46001.59985385615
$ws.Range("C5").Value2 = 46001.60609446955

# Every other row was re-saved in the same bulk write right after, so they
# all share one later timestamp.
$ws.Range("C2").Value2 = 46001.95226143338
$rows = @(6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
foreach ($r in $rows) {
    $ws.Range("C$r").Value2 = 46001.95226143338
}
